# Reorders the semicolon-separated lists in specific cells so the text matches
# the target revision (the underlying data/values are the same multiset of
# tokens, just written back in a different order), and updates the one cell
# whose "Suspect Product Names" value was itself reordered.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "Boredom;Hospitalisation;Aggression;Delinquency"
$ws.Range("E3").Value = "Delinquency;Boredom;Aggression;Hospitalisation"
$ws.Range("E4").Value = "Crime;Delinquency;Aggression;Anger;Imprisonment;Impulsive Behaviour;Drug Withdrawal Syndrome;Libido Decreased;Alcohol Use"
$ws.Range("E5").Value = "Anger;Delinquency"
$ws.Range("E6").Value = "Aggression;Suicide Attempt;Mania;Delinquency;Abnormal Behaviour;Post-Traumatic Stress Disorder;Drug Effect Decreased"
$ws.Range("E7").Value = "Delinquency;Depression;Heart Rate Increased;Intentional Self-Injury;Personality Change;Suicide Attempt;Hand Fracture;Aggression;Hostility;Hypertension"
$ws.Range("E8").Value = "Hallucination, Auditory;Fatigue;Delinquency;Alcohol Interaction;Antisocial Behaviour"
$ws.Range("E9").Value = "Educational Problem;Abnormal Behaviour;Alcohol Use;Drug Abuse;Delinquency;Suicidal Ideation"
$ws.Range("G9").Value = "Required Intervention;Other Outcomes"
$ws.Range("E11").Value = "Delinquency;Schizophrenia"
$ws.Range("B12").Value = "Risperdal Consta;Risperdal"
$ws.Range("E12").Value = "Overdose;Schizophrenia;Depression;Delinquency"
$ws.Range("E13").Value = "Delinquency;Mood Altered;Depression;Condition Aggravated;Anger;Aggression;Social Avoidant Behaviour"
$ws.Range("E14").Value = "Treatment Noncompliance;Suicide Attempt;Delinquency;Personality Disorder;Gambling Disorder"
$ws.Range("Q14").Value = "Detrol;Selegiline Hydrochloride"
$ws.Range("C15").Value = "Oxycodone Hydrochloride;Amitriptyline Hydrochloride;Acetaminophen;Cannabis Sativa Subsp. Sativa Flowering Top;Acetaminophen\Hydrocodone Bitartrate"
$ws.Range("E15").Value = "Alcoholism;Spinal Osteoarthritis;Learning Disability;Male Sexual Dysfunction;Bronchitis;Vomiting;Borderline Mental Impairment;Fall;Delirium;Pollakiuria;Somnolence;Pyrexia;Hand Fracture;Obesity;Pneumoconiosis;Hypertension;Urinary Incontinence;Nocturia;Overdose;Paraesthesia;Asthenia;Inadequate Analgesia;Depressed Level Of Consciousness;Diarrhoea;Generalised Anxiety Disorder;Emotional Distress;Major Depression;Cardiomegaly;Drug Dependence;Bone Pain;Chronic Obstructive Pulmonary Disease;Arthralgia;Toxicity To Various Agents;Hypoaesthesia;Delinquency;Musculoskeletal Stiffness;Speech Disorder;Dyspepsia;Suicidal Ideation;Drug Withdrawal Syndrome;Abnormal Behaviour;Pain In Extremity;Tooth Extraction;Anxiety;Drug Abuser"
$ws.Range("G15").Value = "Other Outcomes;Died;Hospitalized"
$ws.Range("Q15").Value = "Theo-Dur;Phenergan;Lorcet;Ultram;Ambien;Narcan;Flovent;Cogentin;Percocet;Xanax;Haldol;Vioxx;Skelaxin;Buspar;Soma;Baclofen;Pamelor;Tenormin;Toradol;Prednisone;Diazepam;Ativan"
$ws.Range("E16").Value = "Theft;Abnormal Behaviour;Delinquency"
$ws.Range("E17").Value = "Theft;Thinking Abnormal;Legal Problem;Personality Change;Delinquency;Pyromania;Physical Assault"
$ws.Range("E18").Value = "Thinking Abnormal;Paraphilia;Drug Abuser;Theft;Abnormal Behaviour;Obsessive-Compulsive Disorder;Delinquency;Alcoholism;Suicide Attempt;Imprisonment"
$ws.Range("E19").Value = "Depressed Mood;Suicidal Ideation;Physical Assault;Imprisonment;Delinquency;Aggression;Anxiety;Thinking Abnormal"

$wb.Save()
